$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: "Note: Bioinformatics concerns " -> "Note: " + "For our purposes, b"
#         + "ioinformatics concerns " (three runs, same italic Cambria rPr)
# ---------------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("Bioinformatics concerns ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$s1 = $r1.Start
$e1 = $r1.End

# First split "Bioinformatics concerns " away from the preceding "Note: " run
# by nudging formatting (forces the engine to break the run without merging
# neighbouring runs together).
$splitRng1 = $d.Range($s1, $e1)
$splitRng1.Font.Italic = $false
$splitRng1.Font.Italic = $true

# Replace the leading "B" with "For our purposes, b" (lower-cases the "b" and
# inserts the new lead-in text in one content mutation).
$charRng = $d.Range($s1, $s1 + 1)
$charRng.Text = "For our purposes, b"

# Split the freshly-inserted "For our purposes, b" away from the remaining
# "ioinformatics concerns ..." text the same way.
$splitRng2 = $d.Range($s1, $s1 + 19)
$splitRng2.Font.Italic = $false
$splitRng2.Font.Italic = $true

# ---------------------------------------------------------------------------
# Edit 2: append a new sentence after "modelling)." and move the _GoBack
#         bookmark there (it currently sits at the end of the "... research
#         paper in biology is acceptable." paragraph).
# ---------------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("modelling).", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$e2 = $r2.End
$insPoint = $d.Range($e2, $e2)

# Insert the new sentence with a temporary trailing padding character so the
# bookmark insertion point below is never the literal last character of the
# paragraph (collapsed bookmarks placed exactly at a paragraph's final
# position land in the wrong spot), then strip the padding back out.
$insPoint.InsertAfter(" You can also choose to present one of these papers later in the semester. X")

$bmPos = $e2 + 75
$bmRng = $d.Range($bmPos, $bmPos)

$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()
$d.Bookmarks.Add("_GoBack", $bmRng)

# Remove the temporary padding character now that the bookmark is anchored.
$paras = $d.Paragraphs
foreach ($p in $paras) {
    if ($p.Range.Text -like "*You can also choose to present*") {
        $pEnd = $p.Range.End
        $padRng = $d.Range($pEnd - 2, $pEnd - 1)
        $padRng.Text = ""
    }
}

Write-Output "done"
